# Add CU Chat en tiempo real - Iteracion2
# Mirrors the COMUNICACION block already present on ITERACION1 (rows 18-23)
# onto ITERACION2 (rows 19-24), plus a new merged G:I column block holding
# the use-case name "CU Chat en tiempo real" and priority 9.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ITERACION1")
$ws3 = $wb.Worksheets.Item("ITERACION2")

# ---------------------------------------------------------------------
# ITERACION1: the single other change on this sheet - row 21 grew taller
# ---------------------------------------------------------------------
$ws1.Rows.Item(21).RowHeight = 45

# ---------------------------------------------------------------------
# ITERACION2: new requirement rows (19-24) for "CU Chat en tiempo real"
# ---------------------------------------------------------------------

$rows = @(
    @{ r = 19; h = 45; b = "R19"; c = "COMUNICACIÓN"; d = "Pasajero/Conductor"; e = "Accesar fácil al chat"; f = "El sistema permitirá visualizar un botón que facilite el acceso a la opción de entrar al menú del chat." },
    @{ r = 20; h = 30; b = "R20"; c = "COMUNICACIÓN"; d = "Pasajero/Conductor"; e = "Visualizar de chats anteriores"; f = "El sistema permitirá visualizar el contenido de conversaciones anteriores." },
    @{ r = 21; h = 60; b = "R21"; c = "COMUNICACIÓN"; d = "Pasajero/Conductor"; e = "Interactuar en tiempo real"; f = "El sistema debe permitir la interacción de tiempo real, es decir, al momento de escribir un mensaje, que se visualice en el de las otras personas." },
    @{ r = 22; h = 45; b = "R22"; c = "COMUNICACIÓN"; d = "Pasajero/Conductor"; e = "Poder denunciar contenido inapropiado"; f = "El sistema mostrará la opción de poder denunciar un mensaje no apropiado o impertinente para el negocio." },
    @{ r = 23; h = 30; b = "R23"; c = "COMUNICACIÓN"; d = "Pasajero/Conductor"; e = "Notificar cada nuevo mensaje"; f = "Cada nuevo mensaje del pasajero al conductor o viceversa, se notificará" },
    @{ r = 24; h = 30; b = "R24"; c = "COMUNICACIÓN"; d = "Pasajero/Conductor"; e = "Mostrar Mensajes predeterminados"; f = "El sistema mostrará mensajes predeterminados para poder ser de fácil acceso y rápido envío." }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws3.Rows.Item($r).RowHeight = $row.h

    $ws3.Range("B$r").Value = $row.b
    $ws3.Range("C$r").Value = $row.c
    $ws3.Range("D$r").Value = $row.d
    $ws3.Range("E$r").Value = $row.e
    $ws3.Range("F$r").Value = $row.f

    $cells = $ws3.Range("B" + $r + ":F" + $r)
    $cells.WrapText = $true
    $cells.Borders.LineStyle = 1
    $cells.Borders.Weight = 2
}

# New use case name + priority, merged down the whole new block
$ws3.Range("G19").Value = "CU Chat en tiempo real"
$ws3.Range("H19").Value = 9

$gCol = $ws3.Range("G19:G24")
$hCol = $ws3.Range("H19:H24")
$iCol = $ws3.Range("I19:I24")

$gCol.Merge()
$hCol.Merge()
$iCol.Merge()

foreach ($col in @($gCol, $hCol, $iCol)) {
    $col.HorizontalAlignment = -4108  # xlCenter
    $col.VerticalAlignment = -4108    # xlCenter
    $col.WrapText = $true
    $col.Borders.LineStyle = 1
    $col.Borders.Weight = 2
}

# ---------------------------------------------------------------------
# Sheet view / navigation: ITERACION2 becomes the active sheet/tab,
# with the same scroll position + selection captured in the edit.
# ---------------------------------------------------------------------
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

$ws3.Activate()
$ws3.Application.ActiveWindow.ScrollRow = 13
$ws3.Range("F21").Select()
